$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.525.26"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "1.666.12"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4789"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06167"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.48%  "
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("D11").Value = "1.664.63"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5875"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.364"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9999"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "25.522.86"
$ws.Range("E18").Value = "  +2.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006754"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "1.877.07"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.412"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.725"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.274"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "104.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.708"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.967"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07763"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.638"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9994"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04214"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.601"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6097"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.21%  "
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.596"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8619"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01460"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3756"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.841"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.196"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05262"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9993"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.00%  "
